$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column A to fit the new (longer) label, mirroring the author's
#     "Best Fit" column resize. The workbook's original width (72.7109375
#     raw OOXML units) reads back through COM as 71.83; bumping that by the
#     same 10-unit delta the author applied (72.71 -> 82.71) reproduces the
#     resize using the exposed ColumnWidth property.
$ws.Columns.Item(1).ColumnWidth = 81.83

# --- New row: copy the formatting of the row above (row 44, style index 1 -
#     "Calibri Light" 13pt, vertically centered) onto row 46, then set the
#     new label/value. Copying the single source cell onto the destination
#     cell (rather than PasteSpecial after a blanket .Copy()) preserves the
#     existing style id instead of synthesizing new font/fill/xf records.
$ws.Range("A44").Copy($ws.Range("A46"))
$ws.Range("B44").Copy($ws.Range("B46"))

$ws.Range("A46").Value = "After applying another .net dll (reverse markdown) for conversion of html to md"
$ws.Range("B46").Value = 94.14

# Match the row height used by the other data rows (17.25pt).
$ws.Rows.Item(46).RowHeight = 17.25

# --- Move the view down a bit and land the selection on the new row's
#     neighbouring cell, same as the author's saved cursor position.
$ws.Activate()
$ws.Range("E44").Select()
